$d = $word.ActiveDocument

$newText = "Du deltar i en världsomspännande kampanj för att observera och rapportera de svagaste synliga stjärnorna, som ett mått på ljusföroreningarna på orten. Genom att hitta och observera Herkules konstellation på natthimlen kan folk i hela världen lära sig hur belysningen i våra samhällen och omgivningar bidrar till ljusföroreningar. Era bidrag till online-databasen hjälper till att dokumentera den synliga natthimlens över hela världen."

# 1) The first "Kampanjdatum" paragraph has a leading red space run before
#    the "Kampanjdatum för ..." text - match & remove that run too so the
#    whole paragraph collapses into a single, plain (rPr-less) run.
$rngLead = $d.Content.Duplicate
$foundLead = $rngLead.Find.Execute(" Kampanjdatum för Perseus 2018: 30 oktober-8 november och 29 november-8 december", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundLead) {
    $rngLead.Delete()
    $rngLead.InsertAfter($newText)
}

# 2) The remaining 3 "Kampanjdatum" paragraphs (no leading space run).
$rngKamp = $d.Content.Duplicate
while ($rngKamp.Find.Execute("Kampanjdatum för Perseus 2018: 30 oktober-8 november och 29 november-8 december", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rngKamp.Delete()
    $rngKamp.InsertAfter($newText)
    $rngKamp.Collapse(0)
}

# 3) The original descriptive paragraph ("Du deltar ... stjärnbilden (Perseus) ...")
#    gets normalized/collapsed into the same single run text, but mentioning
#    the Hercules constellation instead of Perseus.
$oldDu = "Du deltar i en världsomspännande kampanj för att observera och rapportera de svagaste synliga stjärnorna, som ett mått på ljusföroreningarna på orten. Genom att hitta och observera stjärnbilden (Perseus) på natthimlen kan folk i hela världen lära sig hur belysningen i våra samhällen och omgivningar bidrar till ljusföroreningar. Era bidrag till online-databasen hjälper till att dokumentera den synliga natthimlens över hela världen."
$rngDu = $d.Content.Duplicate
if ($rngDu.Find.Execute($oldDu, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rngDu.Delete()
    $rngDu.InsertAfter($newText)
}
